# Participant130 task_hard.xlsx correction:
# corrected data cleaning for pre/post/total fixation data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the bold/centered/bordered header style from row 1 (A1:AV1),
#    reverting the cells to the default "Normal" style.
$ws.Range("A1:AV1").Style = "Normal"

# 2) Clear the stray "Unnamed: 0" label that used to sit in A1.
$ws.Range("A1").Value = ""

# 3) Correct the fixation-metric values in rows 3-7 (revisit/fixation counts,
#    dwell time ms/%, fixation duration ms) across the affected columns.
$ws.Range("J3").Value = 56
$ws.Range("K3").Value = 13
$ws.Range("L3").Value = 34
$ws.Range("M3").Value = 20
$ws.Range("T3").Value = 13
$ws.Range("U3").Value = 24
$ws.Range("V3").Value = 5
$ws.Range("W3").Value = 21
$ws.Range("AC3").Value = 0
$ws.Range("AD3").Value = 12
$ws.Range("AG3").Value = 22
$ws.Range("AI3").Value = 48
$ws.Range("AR3").Value = 12
$ws.Range("AV3").Value = 4

$ws.Range("J4").Value = 248
$ws.Range("K4").Value = 22
$ws.Range("L4").Value = 130
$ws.Range("M4").Value = 39
$ws.Range("T4").Value = 17
$ws.Range("U4").Value = 41
$ws.Range("V4").Value = 13
$ws.Range("W4").Value = 84
$ws.Range("AC4").Value = 1
$ws.Range("AD4").Value = 21
$ws.Range("AG4").Value = 31
$ws.Range("AI4").Value = 223
$ws.Range("AR4").Value = 15
$ws.Range("AV4").Value = 6

$ws.Range("J5").Value = 88156.2
$ws.Range("K5").Value = 9075.84
$ws.Range("L5").Value = 44600.34
$ws.Range("M5").Value = 16935.12
$ws.Range("R5").Value = 110688.6
$ws.Range("T5").Value = 9193.700000000001
$ws.Range("U5").Value = 15911.19
$ws.Range("V5").Value = 6022.63
$ws.Range("W5").Value = 30674.6
$ws.Range("AC5").Value = 267.14
$ws.Range("AD5").Value = 10143.46
$ws.Range("AG5").Value = 15290.32
$ws.Range("AI5").Value = 81012.67999999999
$ws.Range("AR5").Value = 8642.34
$ws.Range("AV5").Value = 2185.18

$ws.Range("C6").Value = 0.26
$ws.Range("D6").Value = 0.17
$ws.Range("F6").Value = 0.02
$ws.Range("I6").Value = 2.76
$ws.Range("J6").Value = 26.4
$ws.Range("K6").Value = 2.72
$ws.Range("L6").Value = 13.36
$ws.Range("M6").Value = 5.07
$ws.Range("N6").Value = 4.08
$ws.Range("O6").Value = 0.09
$ws.Range("T6").Value = 2.75
$ws.Range("U6").Value = 4.77
$ws.Range("V6").Value = 1.8
$ws.Range("W6").Value = 9.19
$ws.Range("X6").Value = 0.77
$ws.Range("Y6").Value = 0.38
$ws.Range("AA6").Value = 0.09
$ws.Range("AB6").Value = 4.73
$ws.Range("AD6").Value = 3.04
$ws.Range("AE6").Value = 0.3
$ws.Range("AF6").Value = 0.97
$ws.Range("AG6").Value = 4.58
$ws.Range("AH6").Value = 1.86
$ws.Range("AI6").Value = 24.26
$ws.Range("AJ6").Value = 2.79
$ws.Range("AK6").Value = 2.98
$ws.Range("AL6").Value = 0.2
$ws.Range("AO6").Value = 0.97
$ws.Range("AQ6").Value = 0.68
$ws.Range("AR6").Value = 2.59
$ws.Range("AS6").Value = 0.24
$ws.Range("AU6").Value = 0.33
$ws.Range("AV6").Value = 0.65

$ws.Range("J7").Value = 355.47
$ws.Range("K7").Value = 412.54
$ws.Range("L7").Value = 343.08
$ws.Range("M7").Value = 434.23
$ws.Range("T7").Value = 540.8099999999999
$ws.Range("U7").Value = 388.08
$ws.Range("V7").Value = 463.28
$ws.Range("W7").Value = 365.17
$ws.Range("AC7").Value = 267.14
$ws.Range("AD7").Value = 483.02
$ws.Range("AG7").Value = 493.24
$ws.Range("AI7").Value = 363.29
$ws.Range("AR7").Value = 576.16
$ws.Range("AV7").Value = 364.2

# 4) Drop the two trailing blank rows (10 and 11) that were left over from
#    the original export, shrinking the sheet's used range to A1:AV9.
$ws.Rows("10:11").Delete()
